$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Append new key/value rows for level titles (rows 103-110).
# Column A (the keys) is written first, top to bottom, then column B
# (the titles) is written in the same order the original author entered
# them (levels 0-3, then 5, 6, 4, 7) so the shared-string table build
# order matches the source workbook.
$keys = @(
    "level_title_0",
    "level_title_1",
    "level_title_2",
    "level_title_3",
    "level_title_4",
    "level_title_5",
    "level_title_6",
    "level_title_7"
)

for ($i = 0; $i -lt $keys.Length; $i++) {
    $ws.Cells.Item(103 + $i, 1).Value = $keys[$i]
}

$ws.Cells.Item(103, 2).Value = "1 - Polygons #1"
$ws.Cells.Item(104, 2).Value = "2 - Polygons #2"
$ws.Cells.Item(105, 2).Value = "3 - Triangles (Angles)"
$ws.Cells.Item(106, 2).Value = "4 - Triangles (Sides)"
$ws.Cells.Item(108, 2).Value = "6 - Quads (Parallelograms)"
$ws.Cells.Item(109, 2).Value = "7 - Quads (Misc.)"
$ws.Cells.Item(107, 2).Value = "5 - Triangles (All)"
$ws.Cells.Item(110, 2).Value = "8 - Quads (Hierarchy)"

# Update sheet view: selection moves to the new last row, and the view
# scrolls down so the newly added rows are visible.
$ws.Activate() | Out-Null
$ws.Range("B112").Select() | Out-Null

# Reposition the workbook window (best effort - mirrors the recorded
# window move in the source edit).
$excel.ActiveWindow.Left = 9135
$excel.ActiveWindow.Top = 2805
